$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F "想去人数" (want-to-go count) updates
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 1227
$wsExhibition.Range("F4").Value = 20
$wsExhibition.Range("F5").Value = 12496
$wsExhibition.Range("F6").Value = 66
$wsExhibition.Range("F7").Value = 25
$wsExhibition.Range("F8").Value = 26
$wsExhibition.Range("F10").Value = 12378
$wsExhibition.Range("F11").Value = 234
$wsExhibition.Range("F12").Value = 4884
$wsExhibition.Range("F13").Value = 4803
$wsExhibition.Range("F14").Value = 150
$wsExhibition.Range("F15").Value = 71
$wsExhibition.Range("F17").Value = 103
$wsExhibition.Range("F18").Value = 960
$wsExhibition.Range("F19").Value = 9

# Sheet "全部类型" (all types) - same underlying rows, offset by extra entries
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1227
$wsAll.Range("F6").Value = 20
$wsAll.Range("F7").Value = 12496
$wsAll.Range("F8").Value = 66
$wsAll.Range("F9").Value = 25
$wsAll.Range("F10").Value = 26
$wsAll.Range("F12").Value = 12378
$wsAll.Range("F13").Value = 234
$wsAll.Range("F14").Value = 4884
$wsAll.Range("F15").Value = 4803
$wsAll.Range("F16").Value = 150
$wsAll.Range("F17").Value = 71
$wsAll.Range("F19").Value = 103
$wsAll.Range("F20").Value = 960
$wsAll.Range("F21").Value = 9
